$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "ゴシキセイガイインコ。" (row 614) was removed entirely.
# Deleting the row shifts every subsequent row up by one, which matches
# the rest of the diff (it is purely a mechanical renumbering of rows
# 615-674 down to 614-673, plus the dimension shrinking to C673).
$ws.Rows.Item(614).Delete()
